# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect refreshed counts from the data source.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F10").Value = 1007
$wsExhibit.Range("F12").Value = 502
$wsExhibit.Range("F15").Value = 12501

# Sheet "全部类型" (sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F12").Value = 1007
$wsAll.Range("F14").Value = 502
$wsAll.Range("F17").Value = 12501
